{"js": "// Apply the RESUMEN.docx revisions:\n//  1. \"for San Pedro Bakery Company\" -> \"for the company bakery San Pedro\"\n//  2. Move the \"_GoBack\" bookmark so it sits right after \"San Pedro \" (before \"that\")\n//  3. \"store, production, storage, sales)\" -> \"store, production, storage and sales)\"\n\nconst body = context.document.body;\n\n// 1. Re-order \"San Pedro Bakery Company\" -> \"the company bakery San Pedro\"\nlet found = body.search(\"for San Pedro Bakery Company \", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"for the company bakery San Pedro \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. Relocate the \"_GoBack\" bookmark from its old spot (after \"be able to\") to right\n//    before \"that is located\" (i.e. immediately after the newly-placed \"San Pedro \").\nconst existingMark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nexistingMark.load(\"isNullObject\");\nawait context.sync();\nif (!existingMark.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nfound = body.search(\"that is located\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  const startPoint = found.items[0].getRange(Word.RangeLocation.start);\n  startPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// 3. \"store, production, storage, sales)\" -> \"store, production, storage and sales)\"\nfound = body.search(\"store, production, storage, sales)\", { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length > 0) {\n  found.items[0].insertText(\"store, production, storage and sales)\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the RESUMEN.docx revisions:\n#  1. \"for San Pedro Bakery Company\" -> \"for the company bakery San Pedro\"\n#  2. Move the \"_GoBack\" bookmark so it sits right after \"San Pedro \" (before \"that\")\n#  3. \"store, production, storage, sales)\" -> \"store, production, storage and sales)\"\n\n$d = $word.ActiveDocument\n\n# 1. Re-order \"San Pedro Bakery Company\" -> \"the company bakery San Pedro\"\n$d.Content.Find.Execute(\"for San Pedro Bakery Company \", $false, $false, $false, $false, $false, $true, 0, $false, \"for the company bakery San Pedro \", 2) | Out-Null\n\n# 2. Relocate the \"_GoBack\" bookmark from its old spot (after \"be able to\") to right\n#    before \"that is located\" (i.e. immediately after the newly-placed \"San Pedro \").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$r = $d.Content\n$r.Find.Execute(\"that is located\") | Out-Null\n$r.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $r)\n\n# 3. \"store, production, storage, sales)\" -> \"store, production, storage and sales)\"\n$d.Content.Find.Execute(\"store, production, storage, sales)\", $false, $false, $false, $false, $false, $true, 0, $false, \"store, production, storage and sales)\", 2) | Out-Null\n"}
